# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$s1 = New-Object 'object[,]' 11,2
$s1[0,0] = 'Report Generated At'
$s1[0,1] = '2025-05-01 01:54:32 PST'
$s1[1,0] = 'Total URLs Scanned'
$s1[1,1] = 360
$s1[2,0] = 'Total Threats Detected (High/Critical)'
$s1[2,1] = 117
$s1[3,0] = 'Severity - Safe'
$s1[3,1] = 152
$s1[4,0] = 'Severity - Low'
$s1[4,1] = 43
$s1[5,0] = 'Severity - Medium'
$s1[5,1] = 18
$s1[6,0] = 'Severity - High'
$s1[6,1] = 105
$s1[7,0] = 'Severity - Critical'
$s1[7,1] = 12
$s1[8,0] = 'Source - Manual Scans'
$s1[8,1] = 2
$s1[9,0] = 'Source - SMS Scans'
$s1[9,1] = 0
$s1[10,0] = 'Source - Email Scans'
$s1[10,1] = 0
$ws1.Range("A3:B13").Value = $s1

# ---- Sheet 2: Weekly Scans ----
$ws2 = $wb.Worksheets.Item("Weekly Scans")
# Remove the 'Phishing Scans' and 'Safe Scans' columns (C & D)
$ws2.Columns.Item(3).Delete()
$ws2.Columns.Item(3).Delete()
$ws2.Range("B1").Value = "Scans"
$s2 = New-Object 'object[,]' 7,2
$s2[0,0] = 'Fri'
$s2[0,1] = 1
$s2[1,0] = 'Sat'
$s2[1,1] = 8
$s2[2,0] = 'Sun'
$s2[2,1] = 0
$s2[3,0] = 'Mon'
$s2[3,1] = 0
$s2[4,0] = 'Tue'
$s2[4,1] = 3
$s2[5,0] = 'Wed'
$s2[5,1] = 3
$s2[6,0] = 'Thu'
$s2[6,1] = 0
$ws2.Range("A2:B8").Value = $s2

# ---- Sheet 3: Recent Logs ----
$ws3 = $wb.Worksheets.Item("Recent Logs")
$s3 = New-Object 'object[,]' 100,6
$s3[0,0] = '2025-04-30T15:11:07.525000'
$s3[0,1] = 'https://faaarkkbook.com'
$s3[0,2] = 'User Scan'
$s3[0,3] = 'SAFE'
$s3[0,4] = 'Safe Link Verified'
$s3[0,5] = 11.2
$s3[1,0] = '2025-04-30T12:43:31.420000'
$s3[1,1] = 'https://jorogoaol.com'
$s3[1,2] = 'User Scan'
$s3[1,3] = 'SAFE'
$s3[1,4] = 'Safe Link Verified'
$s3[1,5] = 11.2
$s3[2,0] = '2025-04-29T16:12:00.237000'
$s3[2,1] = 'https://jorgol.comm'
$s3[2,2] = 'User Scan'
$s3[2,3] = 'SAFE'
$s3[2,4] = 'Safe Link Verified'
$s3[2,5] = 10.6
$s3[3,0] = '2025-04-29T13:57:24.363000'
$s3[3,1] = 'https://goooooogle.com'
$s3[3,2] = 'User Scan'
$s3[3,3] = 'SAFE'
$s3[3,4] = 'Safe Link Verified'
$s3[3,5] = 11.9
$s3[4,0] = '2025-04-29T08:43:00.350000'
$s3[4,1] = 'https://utoob.com'
$s3[4,2] = 'User Scan'
$s3[4,3] = 'HIGH'
$s3[4,4] = 'Phishing Detected'
$s3[4,5] = 97.3
$s3[5,0] = '2025-04-29T03:57:13.465000'
$s3[5,1] = 'https://www.youtube.com'
$s3[5,2] = 'User Scan'
$s3[5,3] = 'SAFE'
$s3[5,4] = 'Safe Link Verified'
$s3[5,5] = 0
$s3[6,0] = '2025-04-26T04:23:25.830000'
$s3[6,1] = 'https://wwww.facbook.com'
$s3[6,2] = 'User Scan'
$s3[6,3] = 'HIGH'
$s3[6,4] = 'Phishing Detected'
$s3[6,5] = 11.9
$s3[7,0] = '2025-04-26T04:14:37.529000'
$s3[7,1] = 'https://wwww.facbook.com'
$s3[7,2] = 'User Scan'
$s3[7,3] = 'HIGH'
$s3[7,4] = 'Phishing Detected'
$s3[7,5] = 11.9
$s3[8,0] = '2025-04-26T04:11:31.899000'
$s3[8,1] = 'https://wwww.facbook.com'
$s3[8,2] = 'User Scan'
$s3[8,3] = 'HIGH'
$s3[8,4] = 'Phishing Detected'
$s3[8,5] = 11.9
$s3[9,0] = '2025-04-26T03:22:42.027000'
$s3[9,1] = 'https://atendimentocorretora.online'
$s3[9,2] = 'User Scan'
$s3[9,3] = 'HIGH'
$s3[9,4] = 'Phishing Detected'
$s3[9,5] = 96.09999999999999
$s3[10,0] = '2025-04-26T03:20:15.590000'
$s3[10,1] = 'https://atendimentocorretora.online'
$s3[10,2] = 'User Scan'
$s3[10,3] = 'HIGH'
$s3[10,4] = 'Phishing Detected'
$s3[10,5] = 96.09999999999999
$s3[11,0] = '2025-04-26T03:17:48.908000'
$s3[11,1] = 'https://atendimentocorretora.online'
$s3[11,2] = 'User Scan'
$s3[11,3] = 'HIGH'
$s3[11,4] = 'Phishing Detected'
$s3[11,5] = 96.09999999999999
$s3[12,0] = '2025-04-26T03:16:49.344000'
$s3[12,1] = 'https://provasencceja2025.site/'
$s3[12,2] = 'User Scan'
$s3[12,3] = 'SAFE'
$s3[12,4] = 'Safe Link Verified'
$s3[12,5] = 10.6
$s3[13,0] = '2025-04-26T03:16:13.574000'
$s3[13,1] = 'https://www.google.com'
$s3[13,2] = 'User Scan'
$s3[13,3] = 'SAFE'
$s3[13,4] = 'Safe Link Verified'
$s3[13,5] = 0
$s3[14,0] = '2025-04-25T15:56:12.850000'
$s3[14,1] = 'knetww.com/LiveAPP'
$s3[14,2] = 'SMS'
$s3[14,3] = 'HIGH'
$s3[14,4] = 'Phishing Detected'
$s3[14,5] = 0.9
$s3[15,0] = '2025-04-25T15:49:41.113000'
$s3[15,1] = 'https://shopee.ph'
$s3[15,2] = 'User Scan'
$s3[15,3] = 'SAFE'
$s3[15,4] = 'Safe Link Verified'
$s3[15,5] = 0
$s3[16,0] = '2025-04-25T13:01:01.420000'
$s3[16,1] = 'Sugal hindi na need pumunta CAS1NO, pIay online at manalo araw araw, upto 8888 na welcome bonus, 24 hours cash in/out.
Web: win-m.life'
$s3[16,2] = 'SMS'
$s3[16,3] = 'HIGH'
$s3[16,4] = 'Phishing Detected'
$s3[16,5] = 1
$s3[17,0] = '2025-04-25T12:45:13.360000'
$s3[17,1] = 'https://www.youtube.com'
$s3[17,2] = 'User Scan'
$s3[17,3] = 'SAFE'
$s3[17,4] = 'Safe Link Verified'
$s3[17,5] = 0
$s3[18,0] = '2025-04-25T10:57:24.024000'
$s3[18,1] = 'https://www.facebook.com'
$s3[18,2] = 'User Scan'
$s3[18,3] = 'SAFE'
$s3[18,4] = 'Safe Link Verified'
$s3[18,5] = 0
$s3[19,0] = '2025-04-25T08:56:27.561000'
$s3[19,1] = 'https://www.youtube.com'
$s3[19,2] = 'User Scan'
$s3[19,3] = 'SAFE'
$s3[19,4] = 'Safe Link Verified'
$s3[19,5] = 0
$s3[20,0] = '2025-04-25T08:45:07.934000'
$s3[20,1] = 'https://tera.vin/BDOOnline'
$s3[20,2] = 'User Scan'
$s3[20,3] = 'SAFE'
$s3[20,4] = 'Safe Link Verified'
$s3[20,5] = 11.9
$s3[21,0] = '2025-04-25T08:12:03.631000'
$s3[21,1] = 'https://phlpostso-gov.com/ph'
$s3[21,2] = 'User Scan'
$s3[21,3] = 'SAFE'
$s3[21,4] = 'Safe Link Verified'
$s3[21,5] = 11.2
$s3[22,0] = '2025-04-25T08:07:23.586000'
$s3[22,1] = 'https://www.instagram.com'
$s3[22,2] = 'User Scan'
$s3[22,3] = 'SAFE'
$s3[22,4] = 'Safe Link Verified'
$s3[22,5] = 0
$s3[23,0] = '2025-04-25T07:07:47.577000'
$s3[23,1] = 'https://www.facebook.com'
$s3[23,2] = 'User Scan'
$s3[23,3] = 'SAFE'
$s3[23,4] = 'Safe Link Verified'
$s3[23,5] = 0
$s3[24,0] = '2025-04-24T22:45:25.861000'
$s3[24,1] = 'https://lmportal.uc.edu.ph'
$s3[24,2] = 'User Scan'
$s3[24,3] = 'SAFE'
$s3[24,4] = 'Safe Link Verified'
$s3[24,5] = 5.9
$s3[25,0] = '2025-04-24T22:45:24.419000'
$s3[25,1] = 'https://lmportal.uc.edu.ph'
$s3[25,2] = 'User Scan'
$s3[25,3] = 'SAFE'
$s3[25,4] = 'Safe Link Verified'
$s3[25,5] = 5.9
$s3[26,0] = '2025-04-24T22:42:20.723000'
$s3[26,1] = 'https://lmportal.uc.edu.ph'
$s3[26,2] = 'User Scan'
$s3[26,3] = 'SAFE'
$s3[26,4] = 'Safe Link Verified'
$s3[26,5] = 5.9
$s3[27,0] = '2025-04-24T22:39:58.283000'
$s3[27,1] = 'https://www.guayaquil.gob.ec'
$s3[27,2] = 'User Scan'
$s3[27,3] = 'SAFE'
$s3[27,4] = 'Safe Link Verified'
$s3[27,5] = 6.1
$s3[28,0] = '2025-04-24T22:39:55.002000'
$s3[28,1] = 'https://www.guayaquil.gob.ec'
$s3[28,2] = 'User Scan'
$s3[28,3] = 'SAFE'
$s3[28,4] = 'Safe Link Verified'
$s3[28,5] = 6.1
$s3[29,0] = '2025-04-24T22:33:21.311000'
$s3[29,1] = 'https://www.guayaquil.gob.ec'
$s3[29,2] = 'User Scan'
$s3[29,3] = 'SAFE'
$s3[29,4] = 'Safe Link Verified'
$s3[29,5] = 6.1
$s3[30,0] = '2025-04-24T22:33:17.589000'
$s3[30,1] = 'https://www.guayaquil.gob.ec'
$s3[30,2] = 'User Scan'
$s3[30,3] = 'SAFE'
$s3[30,4] = 'Safe Link Verified'
$s3[30,5] = 6.1
$s3[31,0] = '2025-04-24T22:30:11.208000'
$s3[31,1] = 'https://www.guayaquil.gob.ec'
$s3[31,2] = 'User Scan'
$s3[31,3] = 'SAFE'
$s3[31,4] = 'Safe Link Verified'
$s3[31,5] = 6.1
$s3[32,0] = '2025-04-24T22:30:10.807000'
$s3[32,1] = 'https://www.guayaquil.gob.ec'
$s3[32,2] = 'User Scan'
$s3[32,3] = 'SAFE'
$s3[32,4] = 'Safe Link Verified'
$s3[32,5] = 6.1
$s3[33,0] = '2025-04-24T22:18:47.012000'
$s3[33,1] = 'https://www.guayaquil.gob.ec'
$s3[33,2] = 'User Scan'
$s3[33,3] = 'SAFE'
$s3[33,4] = 'Safe Link Verified'
$s3[33,5] = 6.1
$s3[34,0] = '2025-04-24T22:18:45.246000'
$s3[34,1] = 'https://www.guayaquil.gob.ec'
$s3[34,2] = 'User Scan'
$s3[34,3] = 'SAFE'
$s3[34,4] = 'Safe Link Verified'
$s3[34,5] = 6.1
$s3[35,0] = '2025-04-24T22:13:43.584000'
$s3[35,1] = 'https://www.guayaquil.gob.ec'
$s3[35,2] = 'User Scan'
$s3[35,3] = 'SAFE'
$s3[35,4] = 'Safe Link Verified'
$s3[35,5] = 6.1
$s3[36,0] = '2025-04-24T22:13:40.278000'
$s3[36,1] = 'https://www.guayaquil.gob.ec'
$s3[36,2] = 'User Scan'
$s3[36,3] = 'SAFE'
$s3[36,4] = 'Safe Link Verified'
$s3[36,5] = 6.1
$s3[37,0] = '2025-04-24T22:05:21.607000'
$s3[37,1] = 'https://www.guayaquil.gob.ec'
$s3[37,2] = 'User Scan'
$s3[37,3] = 'SAFE'
$s3[37,4] = 'Safe Link Verified'
$s3[37,5] = 6.1
$s3[38,0] = '2025-04-24T22:05:21.576000'
$s3[38,1] = 'https://www.guayaquil.gob.ec'
$s3[38,2] = 'User Scan'
$s3[38,3] = 'SAFE'
$s3[38,4] = 'Safe Link Verified'
$s3[38,5] = 6.1
$s3[39,0] = '2025-04-24T21:43:20.981000'
$s3[39,1] = 'https://www.guayaquil.gob.ec'
$s3[39,2] = 'User Scan'
$s3[39,3] = 'SAFE'
$s3[39,4] = 'Safe Link Verified'
$s3[39,5] = 6.1
$s3[40,0] = '2025-04-24T21:40:53.022000'
$s3[40,1] = 'https://www.guayaquil.gob.ec'
$s3[40,2] = 'User Scan'
$s3[40,3] = 'SAFE'
$s3[40,4] = 'Safe Link Verified'
$s3[40,5] = 6.1
$s3[41,0] = '2025-04-24T21:38:53.424000'
$s3[41,1] = 'https://www.guayaquil.gob.ec'
$s3[41,2] = 'User Scan'
$s3[41,3] = 'SAFE'
$s3[41,4] = 'Safe Link Verified'
$s3[41,5] = 6.1
$s3[42,0] = '2025-04-24T21:36:32.311000'
$s3[42,1] = 'https://teretzurwaalet.webflow.io'
$s3[42,2] = 'User Scan'
$s3[42,3] = 'SAFE'
$s3[42,4] = 'Safe Link Verified'
$s3[42,5] = 2.5
$s3[43,0] = '2025-04-24T21:32:44.583000'
$s3[43,1] = 'https://teretzurwaalet.webflow.io'
$s3[43,2] = 'User Scan'
$s3[43,3] = 'SAFE'
$s3[43,4] = 'Safe Link Verified'
$s3[43,5] = 2.5
$s3[44,0] = '2025-04-24T19:14:22.109000'
$s3[44,1] = 'https://ashleymadisonid.com/'
$s3[44,2] = 'User Scan'
$s3[44,3] = 'HIGH'
$s3[44,4] = 'Phishing Detected'
$s3[44,5] = 96.09999999999999
$s3[45,0] = '2025-04-24T19:13:42.977000'
$s3[45,1] = 'https://grohwtsuppllemments.site/'
$s3[45,2] = 'User Scan'
$s3[45,3] = 'SAFE'
$s3[45,4] = 'Safe Link Verified'
$s3[45,5] = 10.6
$s3[46,0] = '2025-04-24T19:13:14.196000'
$s3[46,1] = 'https://www.tusoatya.online/'
$s3[46,2] = 'User Scan'
$s3[46,3] = 'HIGH'
$s3[46,4] = 'Phishing Detected'
$s3[46,5] = 95.5
$s3[47,0] = '2025-04-24T19:13:14.104000'
$s3[47,1] = 'https://www.tusoatya.online/'
$s3[47,2] = 'User Scan'
$s3[47,3] = 'HIGH'
$s3[47,4] = 'Phishing Detected'
$s3[47,5] = 95.5
$s3[48,0] = '2025-04-24T19:13:13.811000'
$s3[48,1] = 'https://www.tusoatya.online/'
$s3[48,2] = 'User Scan'
$s3[48,3] = 'HIGH'
$s3[48,4] = 'Phishing Detected'
$s3[48,5] = 95.5
$s3[49,0] = '2025-04-24T19:12:38.635000'
$s3[49,1] = 'https://siat.info/D3WYKsWn'
$s3[49,2] = 'User Scan'
$s3[49,3] = 'SAFE'
$s3[49,4] = 'Safe Link Verified'
$s3[49,5] = 16.4
$s3[50,0] = '2025-04-24T19:12:08.528000'
$s3[50,1] = 'https://atendimentocorretora.online'
$s3[50,2] = 'User Scan'
$s3[50,3] = 'HIGH'
$s3[50,4] = 'Phishing Detected'
$s3[50,5] = 96.09999999999999
$s3[51,0] = '2025-04-24T19:11:52.477000'
$s3[51,1] = 'https://central-atendimentoseguro.com'
$s3[51,2] = 'User Scan'
$s3[51,3] = 'SAFE'
$s3[51,4] = 'Safe Link Verified'
$s3[51,5] = 10.6
$s3[52,0] = '2025-04-24T19:11:20.495000'
$s3[52,1] = 'https://gsnews24.com/inst/index.html'
$s3[52,2] = 'User Scan'
$s3[52,3] = 'HIGH'
$s3[52,4] = 'Phishing Detected'
$s3[52,5] = 96.5
$s3[53,0] = '2025-04-24T19:10:57.016000'
$s3[53,1] = 'http://www.personaliteeee.com'
$s3[53,2] = 'User Scan'
$s3[53,3] = 'HIGH'
$s3[53,4] = 'Phishing Detected'
$s3[53,5] = 96.3
$s3[54,0] = '2025-04-24T19:10:33.250000'
$s3[54,1] = 'https://sbi-accessconfirm.zfssw.com/page/'
$s3[54,2] = 'User Scan'
$s3[54,3] = 'SAFE'
$s3[54,4] = 'Safe Link Verified'
$s3[54,5] = 9.5
$s3[55,0] = '2025-04-24T19:10:20.564000'
$s3[55,1] = 'https://avisos-sat.com.mx/'
$s3[55,2] = 'User Scan'
$s3[55,3] = 'SAFE'
$s3[55,4] = 'Safe Link Verified'
$s3[55,5] = 10.6
$s3[56,0] = '2025-04-24T19:09:41.229000'
$s3[56,1] = 'https://avisos-sat.com.mx/'
$s3[56,2] = 'User Scan'
$s3[56,3] = 'SAFE'
$s3[56,4] = 'Safe Link Verified'
$s3[56,5] = 10.6
$s3[57,0] = '2025-04-24T19:09:25.549000'
$s3[57,1] = 'https://biadigitalatendimento.ru.com/'
$s3[57,2] = 'User Scan'
$s3[57,3] = 'SAFE'
$s3[57,4] = 'Safe Link Verified'
$s3[57,5] = 10.6
$s3[58,0] = '2025-04-24T19:09:07.708000'
$s3[58,1] = 'https://validarnetempresax.com/suporte'
$s3[58,2] = 'User Scan'
$s3[58,3] = 'SAFE'
$s3[58,4] = 'Safe Link Verified'
$s3[58,5] = 10.6
$s3[59,0] = '2025-04-24T19:08:42.755000'
$s3[59,1] = 'https://app-usaflex-online.myshopify.com/'
$s3[59,2] = 'User Scan'
$s3[59,3] = 'HIGH'
$s3[59,4] = 'Phishing Detected'
$s3[59,5] = 95.5
$s3[60,0] = '2025-04-24T19:07:12.978000'
$s3[60,1] = 'https://www.pro-bet7k.com/'
$s3[60,2] = 'User Scan'
$s3[60,3] = 'SAFE'
$s3[60,4] = 'Safe Link Verified'
$s3[60,5] = 9.5
$s3[61,0] = '2025-04-24T19:06:55.551000'
$s3[61,1] = 'https://caixageraldepositoseguranca.com/'
$s3[61,2] = 'User Scan'
$s3[61,3] = 'SAFE'
$s3[61,4] = 'Safe Link Verified'
$s3[61,5] = 10.6
$s3[62,0] = '2025-04-24T19:06:31.515000'
$s3[62,1] = 'https://httpss-wwwv-roblox.com'
$s3[62,2] = 'User Scan'
$s3[62,3] = 'HIGH'
$s3[62,4] = 'Phishing Detected'
$s3[62,5] = 96.3
$s3[63,0] = '2025-04-24T19:06:11.094000'
$s3[63,1] = 'https://www-the-graph.xyz'
$s3[63,2] = 'User Scan'
$s3[63,3] = 'SAFE'
$s3[63,4] = 'Safe Link Verified'
$s3[63,5] = 9.5
$s3[64,0] = '2025-04-24T19:05:51.952000'
$s3[64,1] = 'https://robinhood-z.com'
$s3[64,2] = 'User Scan'
$s3[64,3] = 'SAFE'
$s3[64,4] = 'Safe Link Verified'
$s3[64,5] = 8.9
$s3[65,0] = '2025-04-24T19:05:08.333000'
$s3[65,1] = 'https://multipiier-pendle.com'
$s3[65,2] = 'User Scan'
$s3[65,3] = 'HIGH'
$s3[65,4] = 'Phishing Detected'
$s3[65,5] = 96.3
$s3[66,0] = '2025-04-24T19:04:53.651000'
$s3[66,1] = 'https://usualmoney.finance'
$s3[66,2] = 'User Scan'
$s3[66,3] = 'SAFE'
$s3[66,4] = 'Safe Link Verified'
$s3[66,5] = 12.6
$s3[67,0] = '2025-04-24T19:04:29.004000'
$s3[67,1] = 'https://dubai-token2049.cam'
$s3[67,2] = 'User Scan'
$s3[67,3] = 'HIGH'
$s3[67,4] = 'Phishing Detected'
$s3[67,5] = 97.09999999999999
$s3[68,0] = '2025-04-24T19:03:55.231000'
$s3[68,1] = 'http://www.frosttreasuryconnects.com'
$s3[68,2] = 'User Scan'
$s3[68,3] = 'SAFE'
$s3[68,4] = 'Safe Link Verified'
$s3[68,5] = 7.5
$s3[69,0] = '2025-04-24T19:03:25.537000'
$s3[69,1] = 'https://enccejadescomplicado.site/'
$s3[69,2] = 'User Scan'
$s3[69,3] = 'HIGH'
$s3[69,4] = 'Phishing Detected'
$s3[69,5] = 96.09999999999999
$s3[70,0] = '2025-04-24T19:02:53.490000'
$s3[70,1] = 'https://estudeemude.site/'
$s3[70,2] = 'User Scan'
$s3[70,3] = 'HIGH'
$s3[70,4] = 'Phishing Detected'
$s3[70,5] = 96.09999999999999
$s3[71,0] = '2025-04-24T19:02:25.441000'
$s3[71,1] = 'https://provasencceja2025.site/'
$s3[71,2] = 'User Scan'
$s3[71,3] = 'HIGH'
$s3[71,4] = 'Phishing Detected'
$s3[71,5] = 96.90000000000001
$s3[72,0] = '2025-04-24T19:02:00.323000'
$s3[72,1] = 'https://z568zimbra.weebly.com/'
$s3[72,2] = 'User Scan'
$s3[72,3] = 'SAFE'
$s3[72,4] = 'Safe Link Verified'
$s3[72,5] = 0
$s3[73,0] = '2025-04-24T19:01:46.782000'
$s3[73,1] = 'https://bankbac.weebly.com/'
$s3[73,2] = 'User Scan'
$s3[73,3] = 'SAFE'
$s3[73,4] = 'Safe Link Verified'
$s3[73,5] = 0
$s3[74,0] = '2025-04-24T19:01:38.007000'
$s3[74,1] = 'http://secure-paypal-com-login.io'
$s3[74,2] = 'User Scan'
$s3[74,3] = 'SAFE'
$s3[74,4] = 'Safe Link Verified'
$s3[74,5] = 11.2
$s3[75,0] = '2025-04-24T19:01:38.007000'
$s3[75,1] = 'http://secure-paypal-com-login.io'
$s3[75,2] = 'User Scan'
$s3[75,3] = 'SAFE'
$s3[75,4] = 'Safe Link Verified'
$s3[75,5] = 11.2
$s3[76,0] = '2025-04-24T19:01:08.338000'
$s3[76,1] = 'http://apple-support-login-error.ru'
$s3[76,2] = 'User Scan'
$s3[76,3] = 'SAFE'
$s3[76,4] = 'Safe Link Verified'
$s3[76,5] = 11.2
$s3[77,0] = '2025-04-24T19:00:53.705000'
$s3[77,1] = 'http://winner777.pro/sm'
$s3[77,2] = 'User Scan'
$s3[77,3] = 'SAFE'
$s3[77,4] = 'Safe Link Verified'
$s3[77,5] = 17.3
$s3[78,0] = '2025-04-24T19:00:33.763000'
$s3[78,1] = 'http://teretzurwaalet.webflow.io'
$s3[78,2] = 'User Scan'
$s3[78,3] = 'SAFE'
$s3[78,4] = 'Safe Link Verified'
$s3[78,5] = 10.6
$s3[79,0] = '2025-04-24T18:34:45.232000'
$s3[79,1] = 'Kumusta po?'
$s3[79,2] = 'SMS'
$s3[79,3] = 'SAFE'
$s3[79,4] = 'Safe Link Verified'
$s3[79,5] = 0.3
$s3[80,0] = '2025-04-24T18:33:51.681000'
$s3[80,1] = 'Pa may klase ako maya ka na tumawag'
$s3[80,2] = 'SMS'
$s3[80,3] = 'SAFE'
$s3[80,4] = 'Safe Link Verified'
$s3[80,5] = 0
$s3[81,0] = '2025-04-24T18:32:50.262000'
$s3[81,1] = 'Hi, san ka na ba? '
$s3[81,2] = 'SMS'
$s3[81,3] = 'LOW'
$s3[81,4] = 'Safe Link Verified'
$s3[81,5] = 0.6
$s3[82,0] = '2025-04-24T18:32:31.473000'
$s3[82,1] = 'Hi, Kumusta araw mo? '
$s3[82,2] = 'SMS'
$s3[82,3] = 'HIGH'
$s3[82,4] = 'Phishing Detected'
$s3[82,5] = 0.7
$s3[83,0] = '2025-04-24T18:31:09.977000'
$s3[83,1] = 'yLast 3 days, play online BARA-HA,you will have 30% change get XiaoMI/8999P.
web: baraha-p.life'
$s3[83,2] = 'SMS'
$s3[83,3] = 'HIGH'
$s3[83,4] = 'Phishing Detected'
$s3[83,5] = 1
$s3[84,0] = '2025-04-24T18:24:41.336000'
$s3[84,1] = '[Winner777] Free P300 daily login bonus, claim now. https://winner777.pro/sm'
$s3[84,2] = 'SMS'
$s3[84,3] = 'HIGH'
$s3[84,4] = 'Phishing Detected'
$s3[84,5] = 1
$s3[85,0] = '2025-04-24T18:13:51.875000'
$s3[85,1] = 'https://teretzurwaalet.webflow.io'
$s3[85,2] = 'User Scan'
$s3[85,3] = 'SAFE'
$s3[85,4] = 'Safe Link Verified'
$s3[85,5] = 10.6
$s3[86,0] = '2025-04-24T18:08:11.889000'
$s3[86,1] = 'Ang Smart ay nagiimbita para kunin ang inyong 1,888 bonus. Ito ay matatanggap sa inyOng Acc0unt. URL: ww3467.pw/3tZBIuD'
$s3[86,2] = 'SMS'
$s3[86,3] = 'HIGH'
$s3[86,4] = 'Phishing Detected'
$s3[86,5] = 0.9
$s3[87,0] = '2025-04-24T18:00:46.080000'
$s3[87,1] = 'Napaka ganda ng pasok ng Weekend! Dahil dito kay antiviruspedia.net/Live napaka daming B0nu$ naghihintay sayo, hindi lang yan dahil chance mo pang maging VIP!'
$s3[87,2] = 'SMS'
$s3[87,3] = 'HIGH'
$s3[87,4] = 'Phishing Detected'
$s3[87,5] = 0.9
$s3[88,0] = '2025-04-24T17:50:52.032000'
$s3[88,1] = '[''http://wp8luck.com/'']'
$s3[88,2] = 'SMS'
$s3[88,3] = 'HIGH'
$s3[88,4] = 'Phishing Detected'
$s3[88,5] = 1
$s3[89,0] = '2025-04-24T17:35:38.939000'
$s3[89,1] = 'https://landing.twistysnetwork.com/?ats=eyJhIjoxMDc0NCwiYyI6NDQ3NDc0ODgsIm4iOjIsInMiOjYyLCJlIjo4OTQxLCJwIjo1N30='
$s3[89,2] = 'User Scan'
$s3[89,3] = 'HIGH'
$s3[89,4] = 'Phishing Detected'
$s3[89,5] = 97.90000000000001
$s3[90,0] = '2025-04-24T17:26:46.463000'
$s3[90,1] = 'https://www.universityofcebu.net'
$s3[90,2] = 'User Scan'
$s3[90,3] = 'HIGH'
$s3[90,4] = 'Phishing Detected'
$s3[90,5] = 95.8
$s3[91,0] = '2025-04-24T17:26:22.534000'
$s3[91,1] = 'https://www.lmportal.uc.edu.ph'
$s3[91,2] = 'User Scan'
$s3[91,3] = 'SAFE'
$s3[91,4] = 'Safe Link Verified'
$s3[91,5] = 10.6
$s3[92,0] = '2025-04-24T17:22:27.069000'
$s3[92,1] = 'https://gsnews24.com/inst/index.html'
$s3[92,2] = 'User Scan'
$s3[92,3] = 'HIGH'
$s3[92,4] = 'Phishing Detected'
$s3[92,5] = 96.5
$s3[93,0] = '2025-04-24T17:21:31.024000'
$s3[93,1] = 'https://trezor.secure-hardware.io/'
$s3[93,2] = 'User Scan'
$s3[93,3] = 'HIGH'
$s3[93,4] = 'Phishing Detected'
$s3[93,5] = 8.4
$s3[94,0] = '2025-04-24T17:20:55.959000'
$s3[94,1] = 'https://sp-update.info/'
$s3[94,2] = 'User Scan'
$s3[94,3] = 'HIGH'
$s3[94,4] = 'Phishing Detected'
$s3[94,5] = 8
$s3[95,0] = '2025-04-24T17:20:08.418000'
$s3[95,1] = 'https://sbi-accessconfirm.zfssw.com/page/'
$s3[95,2] = 'User Scan'
$s3[95,3] = 'SAFE'
$s3[95,4] = 'Safe Link Verified'
$s3[95,5] = 9.5
$s3[96,0] = '2025-04-24T17:19:38.774000'
$s3[96,1] = 'https://zedny.com.sa/'
$s3[96,2] = 'User Scan'
$s3[96,3] = 'SAFE'
$s3[96,4] = 'Safe Link Verified'
$s3[96,5] = 8.9
$s3[97,0] = '2025-04-24T17:16:56.307000'
$s3[97,1] = 'https://xhhld.com/jcb/'
$s3[97,2] = 'User Scan'
$s3[97,3] = 'SAFE'
$s3[97,4] = 'Safe Link Verified'
$s3[97,5] = 10
$s3[98,0] = '2025-04-24T10:52:19.748000'
$s3[98,1] = 'https://tiktok.com'
$s3[98,2] = 'User Scan'
$s3[98,3] = 'SAFE'
$s3[98,4] = 'Safe Link Verified'
$s3[98,5] = 2.7
$s3[99,0] = '2025-04-24T09:17:35.321000'
$s3[99,1] = 'https://shopee.ph'
$s3[99,2] = 'User Scan'
$s3[99,3] = 'SAFE'
$s3[99,4] = 'Safe Link Verified'
$s3[99,5] = 4.8
$ws3.Range("A2:F101").Value = $s3

